$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45930, 5227.94741537255, 5546.88222686601, 2952, 6781.466379, 172.850049603894),
    @(45931, 4178.04201791952, 5333.3347289325, 3620, 5978.199802, 146.395521375541),
    @(45932, 4215.53182388346, 5873.88668710934, 3620, 6047.636094, 170.249623217745),
    @(45933, 4288.13014228847, 5344.01551986305, 3620, 6184.729748, 150.858963565607),
    @(45934, 1002.96565734096, 3134.46004675596, 3620, 2645.82686, 48.2217187256253),
    @(45935, 843.421923459722, 2381.8093637407, 3620, 2197.488473, 4.82816305337421),
    @(45936, 4367.23135895568, 5228.46999160756, 3620, 6290.900187, 147.172450818828),
    @(45937, 4367.23135895568, 5474.96149315568, 3620, 6290.900187, 157.44293005),
    @(45938, 4367.23135895568, 5187.07436204248, 3620, 6290.900187, 145.447632920283),
    @(45939, 4367.23135895568, 4994.44089594894, 3620, 6290.900187, 137.421238499719),
    @(45940, 5468.88248205865, 4669.66974260547, 3620, 8250.220166, 159.625309439451),
    @(45941, 1519.17029469927, 2733.21719636084, 3620, 3913.743787, 62.8246120275655),
    @(45942, 1415.02104924416, 2751.86297860512, 3620, 3800.437032, 63.2199567233732),
    @(45943, 5868.78419466605, 5418.8574539398, 3620, 9006.878601, 205.706327511406)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowVals[$j]
    }
}
